$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark the date-like text columns (A, D, E) across rows 2-10 as Text
# so Excel does not auto-convert the ISO-formatted strings into date serials.
$ws.Range("A2:A10").NumberFormat = "@"
$ws.Range("D2:D10").NumberFormat = "@"
$ws.Range("E2:E10").NumberFormat = "@"

# Step 2: write the data grid (row 2 is the newly inserted IPO entry; rows 3-10
# carry forward the previous rows 2-9 data, shifted down by one row).
# Row 2
$ws.Range("A2").Value = '2024-04-18'
$ws.Range("B2").Value = '제일엠앤에스'
$ws.Range("C2").Value = 'KB'
$ws.Range("D2").Value = '2024-04-23'
$ws.Range("E2").Value = '2024-04-30'
$ws.Range("F2").Value = 52800000
$ws.Range("G2").Value = 2400000
$ws.Range("H2").Value = '-'
$ws.Range("I2").Value = 15000
$ws.Range("J2").Value = 18000
$ws.Range("K2").Value = '-'
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = '-'
$ws.Range("N2").Value = '-'
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = '-'
$ws.Range("Q2").Value = '-'
$ws.Range("R2").Value = '1438.96 : 1'
$ws.Range("S2").Value = '-'
$ws.Range("T2").Value = '-'
# Row 3
$ws.Range("A3").Value = '2024-04-15'
$ws.Range("B3").Value = '하나33호스팩'
$ws.Range("C3").Value = '하나'
$ws.Range("D3").Value = '2024-04-18'
$ws.Range("E3").Value = '2024-04-24'
$ws.Range("F3").Value = 7000000
$ws.Range("G3").Value = 3500000
$ws.Range("H3").Value = '-'
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = '-'
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = '-'
$ws.Range("N3").Value = '-'
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = '-'
$ws.Range("Q3").Value = '-'
$ws.Range("R3").Value = '2248.41 : 1'
$ws.Range("S3").Value = '-'
$ws.Range("T3").Value = '-'
# Row 4
$ws.Range("A4").Value = '2024-04-11'
$ws.Range("B4").Value = '신한제13호스팩'
$ws.Range("C4").Value = '신한'
$ws.Range("D4").Value = '2024-04-15'
$ws.Range("E4").Value = '2024-04-22'
$ws.Range("F4").Value = 6000000
$ws.Range("G4").Value = 3000000
$ws.Range("H4").Value = '-'
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = '-'
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = '-'
$ws.Range("N4").Value = '-'
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = '-'
$ws.Range("Q4").Value = '-'
$ws.Range("R4").Value = '1337.88 : 1'
$ws.Range("S4").Value = '-'
$ws.Range("T4").Value = '-'
# Row 5
$ws.Range("A5").Value = '2024-04-02'
$ws.Range("B5").Value = '신한제12호스팩'
$ws.Range("C5").Value = '신한'
$ws.Range("D5").Value = '2024-04-05'
$ws.Range("E5").Value = '2024-04-15'
$ws.Range("F5").Value = 10000000
$ws.Range("G5").Value = 5000000
$ws.Range("H5").Value = '-'
$ws.Range("I5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = '-'
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = '-'
$ws.Range("N5").Value = '-'
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = '-'
$ws.Range("Q5").Value = '-'
$ws.Range("R5").Value = '1698.24 : 1'
$ws.Range("S5").Value = '-'
$ws.Range("T5").Value = '-'
# Row 6
$ws.Range("A6").Value = '2024-03-25'
$ws.Range("B6").Value = '아이엠비디엑스'
$ws.Range("C6").Value = '미래'
$ws.Range("D6").Value = '2024-03-28'
$ws.Range("E6").Value = '2024-04-03'
$ws.Range("F6").Value = 32500000
$ws.Range("G6").Value = 2500000
$ws.Range("H6").Value = '-'
$ws.Range("I6").Value = 7700
$ws.Range("J6").Value = 9900
$ws.Range("K6").Value = '-'
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = '-'
$ws.Range("N6").Value = '-'
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = '-'
$ws.Range("Q6").Value = '-'
$ws.Range("R6").Value = '2654.2 : 1'
$ws.Range("S6").Value = '-'
$ws.Range("T6").Value = '-'
# Row 7
$ws.Range("A7").Value = '2024-03-18'
$ws.Range("B7").Value = '하나32호스팩'
$ws.Range("C7").Value = '하나'
$ws.Range("D7").Value = '2024-03-21'
$ws.Range("E7").Value = '2024-03-27'
$ws.Range("F7").Value = 6000000
$ws.Range("G7").Value = 3000000
$ws.Range("H7").Value = '-'
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = '-'
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = '-'
$ws.Range("N7").Value = '-'
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = '-'
$ws.Range("Q7").Value = '-'
$ws.Range("R7").Value = '2389.8 : 1'
$ws.Range("S7").Value = '-'
$ws.Range("T7").Value = '-'
# Row 8
$ws.Range("A8").Value = '2024-03-14'
$ws.Range("B8").Value = '엔젤로보틱스'
$ws.Range("C8").Value = 'NH'
$ws.Range("D8").Value = '2024-03-19'
$ws.Range("E8").Value = '2024-03-26'
$ws.Range("F8").Value = 32000000
$ws.Range("G8").Value = 1600000
$ws.Range("H8").Value = '-'
$ws.Range("I8").Value = 11000
$ws.Range("J8").Value = 15000
$ws.Range("K8").Value = '-'
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = '-'
$ws.Range("N8").Value = '-'
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = '-'
$ws.Range("Q8").Value = '-'
$ws.Range("R8").Value = '2242.016 : 1'
$ws.Range("S8").Value = '-'
$ws.Range("T8").Value = '-'
# Row 9
$ws.Range("A9").Value = '2024-03-12'
$ws.Range("B9").Value = '삼현'
$ws.Range("C9").Value = '한국'
$ws.Range("D9").Value = '2024-03-15'
$ws.Range("E9").Value = '2024-03-21'
$ws.Range("F9").Value = 60000000
$ws.Range("G9").Value = 2000000
$ws.Range("H9").Value = '-'
$ws.Range("I9").Value = 20000
$ws.Range("J9").Value = 25000
$ws.Range("K9").Value = '-'
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = '-'
$ws.Range("N9").Value = '-'
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = '-'
$ws.Range("Q9").Value = '-'
$ws.Range("R9").Value = '1645.13 : 1'
$ws.Range("S9").Value = '-'
$ws.Range("T9").Value = '-'
# Row 10
$ws.Range("A10").Value = '2024-03-04'
$ws.Range("B10").Value = '오상헬스케어'
$ws.Range("C10").Value = 'NH'
$ws.Range("D10").Value = '2024-03-07'
$ws.Range("E10").Value = '2024-03-13'
$ws.Range("F10").Value = 19800000
$ws.Range("G10").Value = 990000
$ws.Range("H10").Value = '-'
$ws.Range("I10").Value = 13000
$ws.Range("J10").Value = 15000
$ws.Range("K10").Value = '-'
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = '-'
$ws.Range("N10").Value = '-'
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = '-'
$ws.Range("Q10").Value = '-'
$ws.Range("R10").Value = '2126.13 : 1'
$ws.Range("S10").Value = '-'
$ws.Range("T10").Value = '-'

# Step 3: restore default (unstyled) cell style for the text columns now that
# the values are safely stored as text, matching the rest of the data rows.
$ws.Range("A2:A10").Style = "Normal"
$ws.Range("D2:D10").Style = "Normal"
$ws.Range("E2:E10").Style = "Normal"
